$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Types" feature row: mark as fully implemented and update the comment to
# reflect that static validation has been added alongside runtime validation.
$ws.Range("B9").Value = "fully implemented"
$ws.Range("C9").Value = "static validation and runtime validation implemented"

# "Cycles" feature row: refine the comment to explain why cycle detection
# still depends on runtime values.
$ws.Range("C10").Value = "Cycles are detected at runtime, they takes into account the current visibilty of the referred Questions which cannot be statically validated since they depends on runtime values"

# Update the saved selection/active cell.
$ws.Range("C11").Select()
